$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing scores for rows 9 and 10 (columns E and F)
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 5
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 5

# Move the active selection down to E11 (frozen pane scrolls accordingly)
$ws.Range("E11").Select()
